# Commit: Add thêm nhân sự Nguyễn Hữu Quang
# This populates the per-person report sheet ("Đơn sale chính") with the
# single service-order row for Lê Thị Ngọc Mi (08-2024) plus its totals
# row, and updates the "Lương" (salary) summary sheet to reflect the
# newly added entry (SÓC TRĂNG branch).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Đơn sale chính"
$ws2 = $wb.Worksheets.Item(2)   # "Lương"

# ---------------------------------------------------------------
# Sheet 1: "Đơn sale chính" - add header row, data row and totals row
# ---------------------------------------------------------------

# Header row (row 1)
$ws1.Range("A1").Value2 = "Tiền tố"
$ws1.Range("B1").Value2 = "Mã dịch vụ"
$ws1.Range("C1").Value2 = "Ngày thực hiện"
$ws1.Range("D1").Value2 = "Cơ sở"
$ws1.Range("E1").Value2 = "Khách hàng"
$ws1.Range("F1").Value2 = "Nguồn khách"
$ws1.Range("G1").Value2 = "Tên dịch vụ"
$ws1.Range("H1").Value2 = "Đơn giá gốc"
$ws1.Range("I1").Value2 = "Sale phụ"
$ws1.Range("J1").Value2 = "Upsale"
$ws1.Range("K1").Value2 = "Đơn giá"
$ws1.Range("L1").Value2 = "Đã thanh toán"
$ws1.Range("M1").Value2 = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("N1").Value2 = "Chiết khấu sale chính"

# Data row (row 2)
$ws1.Range("A2").Value2 = "HD-LUXURY"
$ws1.Range("B2").Value2 = 707
# Force text format so the date-looking string isn't auto-converted to a date serial
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value2 = "08-31-2024"
$ws1.Range("D2").Value2 = "SÓC TRĂNG"
$ws1.Range("E2").Value2 = "dương ngọc hân"
$ws1.Range("F2").Value2 = "Cá nhân"
$ws1.Range("G2").Value2 = "Tiêm Filler"
$ws1.Range("H2").Value2 = 1300000
$ws1.Range("K2").Value2 = 1300000
$ws1.Range("L2").Value2 = 1300000
$ws1.Range("M2").Value2 = 0.1
$ws1.Range("N2").Value2 = 130000

# Totals row (row 3)
$ws1.Range("A3").Value2 = "Tổng"
$ws1.Range("B3").Value2 = 1
$ws1.Range("H3").Value2 = 1300000
$ws1.Range("J3").Value2 = 0
$ws1.Range("K3").Value2 = 1300000
$ws1.Range("L3").Value2 = 1300000
$ws1.Range("M3").Value2 = 0
$ws1.Range("N3").Value2 = 130000

# ---------------------------------------------------------------
# Sheet 2: "Lương" - update the salary summary figures
# ---------------------------------------------------------------

$ws2.Range("B1").Value2 = 5

$ws2.Range("B22").Value2 = 25.5
$ws2.Range("B23").Value2 = 892500
$ws2.Range("B24").Value2 = 2732142.857142857
$ws2.Range("B25").Value2 = 130000

$ws2.Range("B34").Value2 = 3754642.857142857
$ws2.Range("B35").Value2 = 3754642.857142857
